$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, extend the bold/bordered/centered label style (currently on
#     A2:A8) down through the new rows A9:A16 by copying format from an
#     existing styled cell (A7) before any values move around. ---
$ws.Range("A7").Copy()
$ws.Range("A8:A16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 2: Problem 1 -- only F2/H2 change ---
$ws.Range("F2").Value = 65.39
$ws.Range("H2").Value = 65.39

# --- Row 3: Problem 2 ---
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 40
$ws.Range("F3").Value = 213.44
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 3.28

# --- Row 4: Problem 3 (was "Min") ---
$ws.Range("A4").Value = "Problem 3"
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 83.33
$ws.Range("F4").Value = 40.06
$ws.Range("G4").Value = 26.67
$ws.Range("H4").Value = 10.68

# --- Row 5: Problem 4 (was "Max") ---
$ws.Range("A5").Value = "Problem 4"
$ws.Range("C5").Value = 80
$ws.Range("D5").Value = 80
$ws.Range("E5").Value = 75
$ws.Range("F5").Value = 20.83
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 16.67

# --- Row 6: Problem 5 (was "Mean") ---
$ws.Range("A6").Value = "Problem 5"
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 69.44
$ws.Range("F6").Value = 22.32
$ws.Range("G6").Value = 41.67
$ws.Range("H6").Value = 9.82

# --- Row 7: Problem 6 (was "Std") ---
$ws.Range("A7").Value = "Problem 6"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 59.57
$ws.Range("E7").Value = 74
$ws.Range("F7").Value = 21.83
$ws.Range("G7").Value = 44
$ws.Range("H7").Value = 13.93

# --- Row 8: Problem 7 (was "Std / (max - min) %") ---
$ws.Range("A8").Value = "Problem 7"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 83.33
$ws.Range("F8").Value = 10.12
$ws.Range("G8").Value = 83.33
$ws.Range("H8").Value = 10.12

# --- Row 9: Problem 8 (new row) ---
$ws.Range("A9").Value = "Problem 8"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 60
$ws.Range("E9").Value = 79.55
$ws.Range("F9").Value = 15.92
$ws.Range("G9").Value = 47.73
$ws.Range("H9").Value = 9.550000000000001

# --- Row 10: Problem 9 (new row) ---
$ws.Range("A10").Value = "Problem 9"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = 80.36
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 80.36
$ws.Range("H10").Value = 10

# --- Row 11: Problem 10 (new row) ---
$ws.Range("A11").Value = "Problem 10"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 100
$ws.Range("E11").Value = 71.43000000000001
$ws.Range("F11").Value = 19.62
$ws.Range("G11").Value = 71.43000000000001
$ws.Range("H11").Value = 19.62

# --- Row 12: Min (new row) ---
$ws.Range("A12").Value = "Min"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = 3.28

# --- Row 13: Max (new row) ---
$ws.Range("A13").Value = "Max"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = 100
$ws.Range("F13").Value = 213.44
$ws.Range("G13").Value = 100
$ws.Range("H13").Value = 65.39

# --- Row 14: Mean (new row) ---
$ws.Range("A14").Value = "Mean"
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 72
$ws.Range("D14").Value = 72.95699999999999
$ws.Range("E14").Value = 76.644
$ws.Range("F14").Value = 43.953
$ws.Range("G14").Value = 56.51900000000001
$ws.Range("H14").Value = 16.906

# --- Row 15: Std (new row) ---
$ws.Range("A15").Value = "Std"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 28.59681411936962
$ws.Range("D15").Value = 26.71110133342398
$ws.Range("E15").Value = 12.71849501928414
$ws.Range("F15").Value = 61.79540850437208
$ws.Range("G15").Value = 27.73229262871075
$ws.Range("H15").Value = 17.60515088767426

# --- Row 16: Std / (max - min) % (moved down from row 8) ---
$ws.Range("A16").Value = "Std / (max - min) %"
$ws.Range("B16").Value = "inf"
$ws.Range("C16").Value = "inf"
$ws.Range("D16").Value = "inf"
$ws.Range("E16").Value = "inf"
$ws.Range("F16").Value = "inf"
$ws.Range("G16").Value = "inf"
$ws.Range("H16").Value = "inf"
